$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = -8.089900000000002
$ws.Range("E5").Value = 12.24209999999999

$ws.Range("D6").Value = -8.1835
$ws.Range("E6").Value = 12.54810000000001

$ws.Range("C7").Value = -11.14449999999999
$ws.Range("D7").Value = -7.822899999999993

$ws.Range("A8").Value = -20.9462
$ws.Range("D8").Value = -8.0831

$ws.Range("D9").Value = -8.289300000000003

$ws.Range("A10").Value = -20.47659999999998
$ws.Range("D10").Value = -6.361699999999995

$ws.Range("A12").Value = -22.82790000000004
$ws.Range("D12").Value = -8.253499999999997
$ws.Range("E12").Value = 12.70879999999999

$ws.Range("B13").Value = 6.464699999999999

$ws.Range("A18").Value = -22.60450000000003

$ws.Range("E19").Value = 12.95599999999999

$ws.Range("C20").Value = -15.04169999999999
$ws.Range("E20").Value = 12.78969999999999

$ws.Range("E23").Value = 13.20440000000001

$ws.Range("A25").Value = -22.35360000000003
$ws.Range("E25").Value = 13.1817
